$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7287117838859558
$ws.Range("B1").Value = 0.7044810056686401
$ws.Range("C1").Value = 0.6323136687278748
$ws.Range("D1").Value = 2.546707391738892
$ws.Range("E1").Value = 1.491896510124207
